# Insert two new data rows (new rows 69 and 70) into the Ají price sheet,
# pushing the existing rows 69-98 down to 71-100.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A69:A70").EntireRow.Insert()

# New row 69
$ws.Range("A69").Value = 8
$ws.Range("B69").Value = "Terminal La Palmera de La Serena"
$ws.Range("C69").Value = "Coquimbo"
$ws.Range("D69").Value = 44468
$ws.Range("E69").Value = 4
$ws.Range("F69").Value = 100112021
$ws.Range("G69").Value = "Ají"
$ws.Range("H69").Value = "Inferno"
$ws.Range("I69").Value = "Primera"
$ws.Range("J69").Value = 600
$ws.Range("K69").Value = 41000
$ws.Range("L69").Value = 42000
$ws.Range("M69").Value = 41500
$ws.Range("N69").Value = "$/caja 12 kilos"
$ws.Range("O69").Value = "Región de Arica y Parinacota"
$ws.Range("P69").Value = 3458
$ws.Range("Q69").Value = 12
$ws.Range("R69").Value = "Hortaliza"

# New row 70
$ws.Range("A70").Value = 8
$ws.Range("B70").Value = "Terminal La Palmera de La Serena"
$ws.Range("C70").Value = "Coquimbo"
$ws.Range("D70").Value = 44468
$ws.Range("E70").Value = 4
$ws.Range("F70").Value = 100112021
$ws.Range("G70").Value = "Ají"
$ws.Range("H70").Value = "Inferno"
$ws.Range("I70").Value = "Segunda"
$ws.Range("J70").Value = 360
$ws.Range("K70").Value = 35000
$ws.Range("L70").Value = 36000
$ws.Range("M70").Value = 35500
$ws.Range("N70").Value = "$/caja 12 kilos"
$ws.Range("O70").Value = "Región de Arica y Parinacota"
$ws.Range("P70").Value = 2958
$ws.Range("Q70").Value = 12
$ws.Range("R70").Value = "Hortaliza"
